$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "t4 -> t2"
$ws.Range("D2").Value = "t5 -> t2"
$ws.Range("D3").Value = "t6 -> t3"
$ws.Range("D4").Value = "t7 -> t3"
$ws.Range("D5").Value = "t8 -> t4"
$ws.Range("D6").Value = "t3 -> t4"
$ws.Range("D7").Value = "c -> t5"
$ws.Range("D8").Value = "d -> t5"
$ws.Range("D9").Value = "t2, e -> t6"
$ws.Range("D10").Value = "f -> t7"
$ws.Range("D11").Value = "g -> t7"
$ws.Range("D12").Value = "a, b -> t8"
$ws.Range("D13").Value = "t2, t3 -> t1"
